$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065372528896858
$ws.Range("D2").Value = 1.063869726135687
$ws.Range("E2").Value = 1.06945816529658
$ws.Range("F2").Value = 1.077626730693453
$ws.Range("I2").Value = 1.042356725517955
$ws.Range("J2").Value = 1.070327980816843
$ws.Range("K2").Value = 1.066587198548006
$ws.Range("L2").Value = 1.072160619036924
$ws.Range("M2").Value = 1.08030753462386

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066880291374387
$ws.Range("D3").Value = 1.06500833423003
$ws.Range("E3").Value = 1.070778536756297
$ws.Range("F3").Value = 1.078978983542377
$ws.Range("I3").Value = 1.042658868434197
$ws.Range("J3").Value = 1.071488678508453
$ws.Range("K3").Value = 1.067540090210305
$ws.Range("L3").Value = 1.073295914759032
$ws.Range("M3").Value = 1.081476212547099

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067854960537445
$ws.Range("D4").Value = 1.06574402088958
$ws.Range("E4").Value = 1.071632255198699
$ws.Range("F4").Value = 1.079853312558215
$ws.Range("I4").Value = 1.042852545510387
$ws.Range("J4").Value = 1.072238341950822
$ws.Range("K4").Value = 1.068155022228101
$ws.Range("L4").Value = 1.074029328100157
$ws.Range("M4").Value = 1.082231211260303

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.068264489347673
$ws.Range("D5").Value = 1.066053051444675
$ws.Range("E5").Value = 1.071991007378771
$ws.Range("F5").Value = 1.080220724962721
$ws.Range("I5").Value = 1.042933530659424
$ws.Range("J5").Value = 1.072553173250723
$ws.Range("K5").Value = 1.068413147586795
$ws.Range("L5").Value = 1.074337372281104
$ws.Range("M5").Value = 1.082548326332786

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068333238194525
$ws.Range("D6").Value = 1.066104924360006
$ws.Range("E6").Value = 1.072051234723897
$ws.Range("F6").Value = 1.080282406137341
$ws.Range("I6").Value = 1.042947102852417
$ws.Range("J6").Value = 1.072606015753811
$ws.Range("K6").Value = 1.068456465083194
$ws.Range("L6").Value = 1.074389077794368
$ws.Range("M6").Value = 1.082601554678378

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067860433546391
$ws.Range("D7").Value = 1.065748151156951
$ws.Range("E7").Value = 1.071637049450723
$ws.Range("F7").Value = 1.079858222546055
$ws.Range("I7").Value = 1.042853629352201
$ws.Range("J7").Value = 1.072242550022213
$ws.Range("K7").Value = 1.068158472847897
$ws.Range("L7").Value = 1.074033445307505
$ws.Range("M7").Value = 1.082235449687457

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065882285600336
$ws.Range("D8").Value = 1.064254747681384
$ws.Range("E8").Value = 1.069904528414337
$ws.Range("F8").Value = 1.078083872198963
$ws.Range("I8").Value = 1.042459215754636
$ws.Range("J8").Value = 1.070720533643294
$ws.Range("K8").Value = 1.06690957708747
$ws.Range("L8").Value = 1.072544548336917
$ws.Range("M8").Value = 1.080702748746967

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06238894796998
$ws.Range("D9").Value = 1.061614823044062
$ws.Range("E9").Value = 1.066846418508554
$ws.Range("F9").Value = 1.074951907571962
$ws.Range("I9").Value = 1.041750136788558
$ws.Range("J9").Value = 1.068027732235108
$ws.Range("K9").Value = 1.06469604531473
$ws.Range("L9").Value = 1.069911549159661
$ws.Range("M9").Value = 1.077992444153969

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060054558507436
$ws.Range("D10").Value = 1.059849014048126
$ws.Range("E10").Value = 1.064803900057788
$ws.Range("F10").Value = 1.072860051467304
$ws.Range("I10").Value = 1.041267872235626
$ws.Range("J10").Value = 1.066224975744576
$ws.Range("K10").Value = 1.063211512394594
$ws.Range("L10").Value = 1.068149646741638
$ws.Range("M10").Value = 1.076178924115093

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059042349516443
$ws.Range("D11").Value = 1.059082955269732
$ws.Range("E11").Value = 1.063918505300229
$ws.Range("F11").Value = 1.071953268280647
$ws.Range("I11").Value = 1.041056763240856
$ws.Range("J11").Value = 1.065442506219419
$ws.Range("K11").Value = 1.062566544148674
$ws.Range("L11").Value = 1.067385105036826
$ws.Range("M11").Value = 1.075392011819893

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.058666151541671
$ws.Range("D12").Value = 1.058798184041609
$ws.Range("E12").Value = 1.063589478687136
$ws.Range("F12").Value = 1.071616293249205
$ws.Range("I12").Value = 1.040978002916814
$ws.Range("J12").Value = 1.065151576543331
$ws.Range("K12").Value = 1.062326645857332
$ws.Range("L12").Value = 1.067100870590837
$ws.Range("M12").Value = 1.075099464919253

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058746857295587
$ws.Range("D13").Value = 1.058859278599848
$ws.Range("E13").Value = 1.063660062964154
$ws.Range("F13").Value = 1.071688582662067
$ws.Range("I13").Value = 1.040994912915216
$ws.Range("J13").Value = 1.065213994963236
$ws.Range("K13").Value = 1.062378119787077
$ws.Range("L13").Value = 1.067161851225112
$ws.Range("M13").Value = 1.075162228740697

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059011257376326
$ws.Range("D14").Value = 1.059059420565754
$ws.Range("E14").Value = 1.063891310983164
$ws.Range("F14").Value = 1.071925417022796
$ws.Range("I14").Value = 1.041050259936621
$ws.Range("J14").Value = 1.065418463741438
$ws.Range("K14").Value = 1.062546720805812
$ws.Range("L14").Value = 1.067361615251592
$ws.Range("M14").Value = 1.075367834986453

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.059174133809282
$ws.Range("D15").Value = 1.059182704958273
$ws.Range("E15").Value = 1.064033770262268
$ws.Range("F15").Value = 1.072071317725525
$ws.Range("I15").Value = 1.041084315290195
$ws.Range("J15").Value = 1.065544405683791
$ws.Range("K15").Value = 1.062650557810969
$ws.Range("L15").Value = 1.067484663226367
$ws.Range("M15").Value = 1.075494482143668

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.060121705269048
$ws.Range("D16").Value = 1.059899823851725
$ws.Range("E16").Value = 1.064862639850438
$ws.Range("F16").Value = 1.072920210209623
$ws.Range("I16").Value = 1.041281834550928
$ws.Range("J16").Value = 1.066276865919066
$ws.Range("K16").Value = 1.06325427101164
$ws.Range("L16").Value = 1.068200352142035
$ws.Range("M16").Value = 1.076231113756611

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.060715710773628
$ws.Range("D17").Value = 1.060349261535749
$ws.Range("E17").Value = 1.065382303904696
$ws.Range("F17").Value = 1.073452427387221
$ws.Range("I17").Value = 1.041405120132821
$ws.Range("J17").Value = 1.066735815867277
$ws.Range("K17").Value = 1.063632384045327
$ws.Range("L17").Value = 1.06864884576039
$ws.Range("M17").Value = 1.076692738838734

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061062049294344
$ws.Range("D18").Value = 1.060611271177061
$ws.Range("E18").Value = 1.065685321878091
$ws.Range("F18").Value = 1.073762764986951
$ws.Range("I18").Value = 1.041476810148868
$ws.Range("J18").Value = 1.067003334132588
$ws.Range("K18").Value = 1.063852723293633
$ws.Range("L18").Value = 1.068910287941771
$ws.Range("M18").Value = 1.076961838471457

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061180119140652
$ws.Range("D19").Value = 1.060700586132151
$ws.Range("E19").Value = 1.065788627527204
$ws.Range("F19").Value = 1.073868566034593
$ws.Range("I19").Value = 1.041501217263119
$ws.Range("J19").Value = 1.067094520674275
$ws.Range("K19").Value = 1.063927818191211
$ws.Range("L19").Value = 1.068999406573157
$ws.Range("M19").Value = 1.077053567773304

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.060651993607725
$ws.Range("D20").Value = 1.060301055615952
$ws.Range("E20").Value = 1.065326558559462
$ws.Range("F20").Value = 1.073395335459162
$ws.Range("I20").Value = 1.041391915563812
$ws.Range("J20").Value = 1.066686593485869
$ws.Range("K20").Value = 1.063591837609974
$ws.Range("L20").Value = 1.068600742871539
$ws.Range("M20").Value = 1.076643227265906

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058933404236012
$ws.Range("D21").Value = 1.059000489936498
$ws.Range("E21").Value = 1.063823218435758
$ws.Range("F21").Value = 1.07185567954249
$ws.Range("I21").Value = 1.041033971156993
$ws.Range("J21").Value = 1.065358260713144
$ws.Range("K21").Value = 1.062497081108907
$ws.Range("L21").Value = 1.067302796657958
$ws.Range("M21").Value = 1.075307296088613

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.05785159093597
$ws.Range("D22").Value = 1.05818148072005
$ws.Range("E22").Value = 1.062877128984109
$ws.Range("F22").Value = 1.070886734873252
$ws.Range("I22").Value = 1.040806920128906
$ws.Range("J22").Value = 1.06452142966494
$ws.Range("K22").Value = 1.061806862160726
$ws.Range("L22").Value = 1.066485279519102
$ws.Range("M22").Value = 1.074465878792817

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.058425203407647
$ws.Range("D23").Value = 1.058615777064443
$ws.Range("E23").Value = 1.063378754224119
$ws.Range("F23").Value = 1.071400478206472
$ws.Range("I23").Value = 1.040927474057438
$ws.Range("J23").Value = 1.06496520843903
$ws.Range("K23").Value = 1.062172942089881
$ws.Range("L23").Value = 1.066918799736614
$ws.Range("M23").Value = 1.074912070609507

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.060680785069293
$ws.Range("D24").Value = 1.060322838228284
$ws.Range("E24").Value = 1.065351747770016
$ws.Range("F24").Value = 1.073421133139955
$ws.Range("I24").Value = 1.041397882821241
$ws.Range("J24").Value = 1.066708835516142
$ws.Range("K24").Value = 1.063610159440827
$ws.Range("L24").Value = 1.068622478980505
$ws.Range("M24").Value = 1.076665599904441

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.06329299933894
$ws.Range("D25").Value = 1.062298321835475
$ws.Range("E25").Value = 1.067637657734829
$ws.Range("F25").Value = 1.075762257629974
$ws.Range("I25").Value = 1.041935126903195
$ws.Range("J25").Value = 1.068725196677537
$ws.Range("K25").Value = 1.065269838600204
$ws.Range("L25").Value = 1.070593380743781
$ws.Range("M25").Value = 1.078694275121936
